$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.996.40'
$ws.Range("E2").Value = '  -2.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.126.76'
$ws.Range("E3").Value = '  -5.60%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.13'
$ws.Range("E5").Value = '  -2.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.61'
$ws.Range("E6").Value = '  -5.50%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.123.74'
$ws.Range("E8").Value = '  -5.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.506'
$ws.Range("E9").Value = '  -2.67%  '
$ws.Range("E10").Value = '  -6.25%  '
$ws.Range("E11").Value = '  -4.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -4.14%  '
$ws.Range("E13").Value = '  -6.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.86'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.638.11'
$ws.Range("E15").Value = '  -5.70%  '
$ws.Range("E16").Value = '  -2.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.125.09'
$ws.Range("E17").Value = '  -5.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.019.00'
$ws.Range("E18").Value = '  -2.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.50'
$ws.Range("E19").Value = '  -5.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '449.90'
$ws.Range("E20").Value = '  -6.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.77'
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.699'
$ws.Range("E22").Value = '  -4.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("E23").Value = '  -6.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.22'
$ws.Range("E24").Value = '  -3.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.87'
$ws.Range("E25").Value = '  -1.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  -3.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.72'
$ws.Range("E29").Value = '  -8.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.62'
$ws.Range("E30").Value = '  -6.26%  '
$ws.Range("E31").Value = '  -8.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.03'
$ws.Range("E32").Value = '  -6.13%  '
$ws.Range("E33").Value = '  -4.36%  '
$ws.Range("E34").Value = '  -8.24%  '
$ws.Range("E35").Value = '  -8.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.76'
$ws.Range("E36").Value = '  -4.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '50.72'
$ws.Range("E37").Value = '  -4.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0688'
$ws.Range("E38").Value = '  -7.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0383'
$ws.Range("E39").Value = '  -4.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.63'
$ws.Range("E40").Value = '  -4.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.99'
$ws.Range("E41").Value = '  -4.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '392.38'
$ws.Range("E42").Value = '  -9.61%  '
$ws.Range("E43").Value = '  -3.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.747.84'
$ws.Range("E44").Value = '  -10.44%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.248'
$ws.Range("E45").Value = '  -6.33%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  -4.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.33'
$ws.Range("E48").Value = '  -1.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.04'
$ws.Range("E49").Value = '  -5.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.12'
$ws.Range("E50").Value = '  -5.94%  '
$ws.Range("E51").Value = '  -3.78%  '
